$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B (Buying Opportunity) for rows 2-9; rows 10-13 lose their value ---
$ws.Range("B2").Value = "NSE:EIMCOELECO"
$ws.Range("B3").Value = "NSE:FAIRCHEMOR"
$ws.Range("B4").Value = "NSE:JASH"
$ws.Range("B5").Value = "NSE:KIRIINDUS"
$ws.Range("B6").Value = "NSE:METROBRAND"
$ws.Range("B7").Value = "NSE:NSLNISP"
$ws.Range("B8").Value = "NSE:OIL"
$ws.Range("B9").Value = "NSE:PALREDTEC"
$ws.Range("B10:B13").ClearContents()

# --- Update column C (support Zone) for rows 2-13 ---
$ws.Range("C2").Value = "NSE:21STCENMGM"
$ws.Range("C3").Value = "NSE:A2ZINFRA"
$ws.Range("C4").Value = "NSE:DEEPAKNTR"
$ws.Range("C5").Value = "NSE:GOLDBEES"
$ws.Range("C6").Value = "NSE:HDFCGOLD"
$ws.Range("C7").Value = "NSE:HDFCSILVER"
$ws.Range("C8").Value = "NSE:ITBEES"
$ws.Range("C9").Value = "NSE:IVZINGOLD"
$ws.Range("C10").Value = "NSE:ONEPOINT"
$ws.Range("C11").Value = "NSE:PANACEABIO"
$ws.Range("C12").Value = "NSE:QGOLDHALF"
$ws.Range("C13").Value = "NSE:RVNL"

# --- Columns D, E, F (long buildup / Short buildup / FII ENTERING) are now empty for rows 2-13 ---
$ws.Range("D2:F13").ClearContents()

# --- Rows 14-17 (ranks 12-15) are removed entirely, shrinking the used range to A1:F13 ---
$ws.Rows("14:17").Delete()
